# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-09 (serial 45178) to 2023-09-10 (serial 45179).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 151; $row++) {
    $ws.Cells.Item($row, 3).Value = 45179
}
